$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "Absent" (column H) to 1 for rows 3-10
for ($r = 3; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Row 11: "Total Attendance Count" (D) and "Real" (E) become 1
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 1

# Set "Absent" (column H) to 1 for rows 12-15
for ($r = 12; $r -le 15; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Row 16: "Invalid" (G) and "Absent" (H) become 1
$ws.Cells.Item(16, 7).Value = 1
$ws.Cells.Item(16, 8).Value = 1

# Set "Absent" (column H) to 1 for rows 17-18
for ($r = 17; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
